$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $range = $ws.Range($cell)
    $origStyle = $range.Style
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = $origStyle
}

Set-TextValue 'D2' '28.402.12'
Set-TextValue 'E2' '  +4.28%  '
Set-TextValue 'D3' '1.596.12'
Set-TextValue 'E3' '  +2.02%  '
Set-TextValue 'E4' '  -0.10%  '
Set-TextValue 'D5' '214.75'
Set-TextValue 'E5' '  +1.99%  '
Set-TextValue 'E6' '  +1.20%  '
Set-TextValue 'D8' '24.05'
Set-TextValue 'E8' '  +8.52%  '
Set-TextValue 'E9' '  +0.78%  '
Set-TextValue 'D10' '0.0602'
Set-TextValue 'E10' '  +0.88%  '
Set-TextValue 'E11' '  +1.98%  '
Set-TextValue 'D12' '1.823.34'
Set-TextValue 'E12' '  +2.01%  '
Set-TextValue 'D13' '1.589.16'
Set-TextValue 'E13' '  +1.45%  '
Set-TextValue 'E14' '  +2.79%  '
Set-TextValue 'E15' '  -0.05%  '
Set-TextValue 'D16' '28.403.88'
Set-TextValue 'E16' '  +4.39%  '
Set-TextValue 'D17' '63.19'
Set-TextValue 'E17' '  +2.00%  '
Set-TextValue 'D18' '228.41'
Set-TextValue 'E18' '  +4.66%  '
Set-TextValue 'D19' '0.0₃0712'
Set-TextValue 'E19' '  +1.57%  '
Set-TextValue 'E20' '  +0.68%  '
Set-TextValue 'E21' '  -0.10%  '
Set-TextValue 'D22' '4.11'
Set-TextValue 'E22' '  -0.78%  '
Set-TextValue 'E23' '  -0.27%  '
Set-TextValue 'E24' '  +0.63%  '
Set-TextValue 'D25' '151.96'
Set-TextValue 'E25' '  +0.40%  '
Set-TextValue 'E26' '  +1.25%  '
Set-TextValue 'E27' '  +0.63%  '
Set-TextValue 'E28' '  -0.46%  '
Set-TextValue 'E29' '  -0.01%  '
Set-TextValue 'E30' '  +0.93%  '
Set-TextValue 'E31' '  +1.21%  '
Set-TextValue 'E32' '  +0.12%  '
Set-TextValue 'D33' '3.15'
Set-TextValue 'E33' '  -0.44%  '
Set-TextValue 'D34' '1.399.01'
Set-TextValue 'E34' '  -4.10%  '
Set-TextValue 'E35' '  -1.25%  '
Set-TextValue 'E36' '  -5.31%  '
Set-TextValue 'E37' '  +0.34%  '
Set-TextValue 'E38' '  +0.90%  '
Set-TextValue 'E39' '  +7.63%  '
Set-TextValue 'D40' '0.542'
Set-TextValue 'E40' '  +0.27%  '
Set-TextValue 'E41' '  +0.24%  '
Set-TextValue 'D42' '5.74'
Set-TextValue 'E42' '  -2.39%  '
Set-TextValue 'E43' '  -0.12%  '
Set-TextValue 'E44' '  +8.01%  '
Set-TextValue 'D45' '0.983'
Set-TextValue 'E45' '  -0.43%  '
Set-TextValue 'D46' '64.52'
Set-TextValue 'E46' '  +0.08%  '
Set-TextValue 'D47' '1.733.09'
Set-TextValue 'E47' '  +1.94%  '
Set-TextValue 'D48' '87.63'
Set-TextValue 'E48' '  +2.05%  '
Set-TextValue 'E49' '  +0.01%  '
Set-TextValue 'E50' '  -1.55%  '
Set-TextValue 'E51' '  +0.10%  '
